$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "index" (first sheet)
# ---------------------------------------------------------------------
$wsIndex = $wb.Worksheets.Item("index")

# The long "広大moodleとは" HTML blurb that used to live in B6 is removed;
# the cell goes back to being empty (keeping its existing B-column style).
$wsIndex.Range("B6").ClearContents()
# Row height was an explicit 243pt to fit the long text - now that the
# text is gone, let the row return to the sheet's default height.
$wsIndex.Rows.Item(6).AutoFit()

# The cursor/selection on this sheet moves off this sheet entirely (see
# below) - but its stored selection position still changes from B14 to B12.
$wsIndex.Activate()
$wsIndex.Range("B12").Select()

# ---------------------------------------------------------------------
# Sheet "p1" (second sheet)
# ---------------------------------------------------------------------
$wsP1 = $wb.Worksheets.Item("p1")

# Row 6 ("本書では、広大 moodle の基本的な使い方を簡単にご紹介します。...")
# is removed outright, shifting everything below it up by one row.
$wsP1.Rows.Item(6).Delete()

# B2 changes from the generic "はじめに" header text to the page's own
# title text.
$wsP1.Range("B2").Value = "広大moodleとは"

# The row that used to hold the "<h3><a name=""moodle""></a>..." anchor
# heading (now row 6, after the deletion above) is cleared out - the
# following row already carries the full paragraph text that used to sit
# further down, so this becomes a blank spacer row again.
$wsP1.Range("B6:D6").Clear()

# B5 picks up the same (wrap-text) formatting as the title cell above it,
# even though it stays empty.
$wsP1.Range("B4").Copy()
$wsP1.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# This sheet becomes the active tab with the selection sitting on B7.
$wsP1.Activate()
$wsP1.Range("B7").Select()
